$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix row 5 (existing data bug fix)
$ws.Range("B5").Value = 40
$ws.Range("D5").Value = "Georgia"
$ws.Range("E5").Value = "Acme"
$ws.Range("G5").Value = "Consult #1"

# Add new row 7
$ws.Range("A7").Value = "b@gmail.com"
$ws.Range("B7").Value = 40
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = "Georgia"
$ws.Range("E7").Value = "Acme"
$ws.Range("F7").Value = "09/14–09/27"
$ws.Range("G7").Value = "Special Consult #2"
